$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the erroring S21 cell (column S has no source data, so AVERAGE(S2:S16) -> #DIV/0!).
# Clearing its formula/content removes the error while keeping its existing style.
$ws.Range("S21").ClearContents()

# Add a new summary row 25: average of the second group of students (rows 8-16)
# across columns O:X, mirroring the percentage style already used for O21:X21 (style 4).
# Column S is intentionally left blank (no data / no formula), matching S21/S8..S16.
$ws.Range("O25").Formula = "=AVERAGE(O8:O16)"
$ws.Range("P25:X25").FormulaR1C1 = "=AVERAGE(R[-17]C:R[-9]C)"
$ws.Range("S25").ClearContents()

# Apply the percentage style (same style used by O21:X21) across the whole O25:X25 band,
# including the blank S25 cell, so every cell in the row shares a consistent format.
$ws.Range("O25:X25").Style = "Percent"
$ws.Range("O25:X25").Font.Name = $ws.Range("O21").Font.Name
$ws.Range("O25:X25").Font.Size = $ws.Range("O21").Font.Size

# Update the active selection to match where the user ended up after this edit.
$ws.Range("S25").Select() | Out-Null
